# Applies the "Git Fix, Inbetriebnahme Navigation Stack" edit:
#  1. Merge the two runs of the "Der Vorteil ist..." paragraph into one run
#     (no textual change - just a harmless Find/Replace that spans the
#     run boundary, which causes identically-formatted adjacent runs to
#     coalesce, matching the target OOXML).
#  2. Merge the two runs of the "Daraus wuerde ich..." / "Das heisst fuer
#     morgen..." paragraph into one run the same way (again no textual
#     change at this point).
#  3. Insert a new "28.11.2017:" heading paragraph (bold, underlined,
#     matching the style of the existing "27.11.2017:" heading) right
#     after the "Daraus wuerde ich..." paragraph.
#  4. Insert a new body paragraph after that heading with the first new
#     diary entry (this is what textually used to be "Das heisst fuer
#     morgen..." - it now lives in its own paragraph with new wording).
#  5. Insert the second new diary sentence into the (until now empty)
#     run that precedes the manual page break.

$d = $word.ActiveDocument

function Get-ParaIndexAtRangeStart($rng) {
    $tmp = $d.Content
    $tmp.End = $rng.Start
    return $tmp.Paragraphs.Count
}

# --- 1) Force-merge the two runs in the "Der Vorteil ist..." paragraph ---
$rng1 = $d.Content
$rng1.Find.Execute("An der Stelle muss man aber auch erwähnen", $true, $false, $false, $false, $false, $true, 1, $false, "An der Stelle muss man aber auch erwähnen", 2) | Out-Null

# --- 2) Force-merge the two runs in the "Daraus würde ich..." paragraph ---
$rng2 = $d.Content
$rng2.Find.Execute("für die Realität zu definieren. Das heißt für morgen", $true, $false, $false, $false, $false, $true, 1, $false, "für die Realität zu definieren. Das heißt für morgen", 2) | Out-Null

# --- Locate the (now single-run) "Daraus würde ich..." paragraph ---
$anchorRng = $d.Content
$anchorRng.Find.Execute("Navigationsstack soll sowohl mit als auch ohne AMCL laufen.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$anchorIdx = Get-ParaIndexAtRangeStart($anchorRng)

# --- 3) Insert the new "28.11.2017:" heading paragraph right after it ---
$anchorPara = $d.Paragraphs.Item($anchorIdx)
$hRng = $anchorPara.Range
$hRng.Collapse(0)
$hRng.InsertParagraphAfter()

$headingIdx = $anchorIdx + 1
$headingPara = $d.Paragraphs.Item($headingIdx)
$headingPara.Range.Font.Bold = $true
$headingPara.Range.Font.BoldBi = $true
$headingPara.Range.Font.Underline = 1
$headingPara.Range.Text = "28.11.2017:"

# --- 4) Insert the first new body paragraph right after the heading ---
$headingPara2 = $d.Paragraphs.Item($headingIdx)
$bRng = $headingPara2.Range
$bRng.Collapse(0)
$bRng.InsertParagraphAfter()

$body1Idx = $headingIdx + 1
$body1Para = $d.Paragraphs.Item($body1Idx)
$body1Para.Range.Font.Bold = $false
$body1Para.Range.Font.BoldBi = $false
$body1Para.Range.Font.Underline = 0
$body1Para.Range.Text = "Die ersten Aufgaben heute bestand darin den Navigation-Stack in der Standard-Konfiguration zum Laufen zu bringen. Als erster Punkt stand die Erstellung der Korridor-Karte auf dem Plan. Die Lösung von dem Problem bringt MATLAB und der map_server. Die bei der Kartographierung erstelle Karte wird über eine MATLAB-ROS-Node publiziert, von dem map_server empfangen, der die Karte wiederum abspeichert."

# --- 5) Insert the second new sentence before the manual page break ---
# (this is the paragraph that immediately follows the body paragraph we
# just wrote; it contains only an empty run followed by the page-break run)
$pageBreakIdx = $body1Idx + 1
$pbPara = $d.Paragraphs.Item($pageBreakIdx)
$insRng = $pbPara.Range
$insRng.Collapse(1)
$insRng.Text = "Als nächstes war dann der Navigation-Stack dran, wobei das Ganze sowohl mit als auch ohne AMC-Lokalisierung funktionieren sollte. Außerdem sollen playground und Korridor-Karte verwendet werden können. Nach Tutorial hat der Standard-Fall mit Playground und AMCL recht schnell funktioniert, nach paar Problemchen auch ohne AMCL (hier mus sman manuell die Transformation von Odometrie zu Karten-Frame vorgeben, was ohne Lokalisierung die Identitätsabbildung ist). Um den Korridor zum Laufenzu bringen war der map_server schwierig, weil in dem gespeicherten yaml-file noch eine NaN-Wert war, der den Server zum Absturz bringt. Den Faller kann man allerdings händisch korrigieren, woraufhin das Ganz funktoiniet."

Write-Output "done"
